$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct/reorder the "Causas" labels and figures (row 5-12) ---
$ws.Range("B5").Value = "Edema, proteinuria y trastornos hipertensivos"
$ws.Range("C5").Value = 14
$ws.Range("D5").Value = 5.6

$ws.Range("B6").Value = "Hemorragia en el embarazo, parto y el puerperio"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0.4

$ws.Range("B7").Value = "Complicaciones predominantes"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 0.8

$ws.Range("B8").Value = "Sepsis puerperal y otras infecciones"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0

$ws.Range("B9").Value = "Complicaciones relacionadas con el puerperio"
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 1.2

$ws.Range("B10").Value = "Muertes obstétricas de causa no específica"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0

$ws.Range("B11").Value = "Cualquier causa obstétrica que ocurre después de 42 días pero antes de un año del parto"
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 2.4

$ws.Range("B12").Value = "Muertes obstétricas indirectas"
$ws.Range("C12").Value = 12
$ws.Range("D12").Value = 4.8

# --- Widen column B so the longer labels are readable ---
$ws.Columns.Item(2).ColumnWidth = 59

# --- Move the active selection like the author left it ---
$ws.Range("F12").Select() | Out-Null

# --- Allow iterative calculation (author turned this on) ---
$excel.Iteration = $true | Out-Null
